$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for cells whose numeric-looking values must keep exact
# string representation (trailing zeros / leading zeros) instead of being
# auto-converted to a number by Excel.
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "66.417.43"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "3.839.94"
$ws.Range("E3").Value = "  +8.88%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "426.81"
$ws.Range("E5").Value = "  +8.73%  "
$ws.Range("D6").Value = "131.72"
$ws.Range("E6").Value = "  +7.58%  "
$ws.Range("D7").Value = "3.833.69"
$ws.Range("E7").Value = "  +8.95%  "
$ws.Range("D8").Value = "0.614"
$ws.Range("E8").Value = "  +4.61%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "0.731"
$ws.Range("E10").Value = "  +8.24%  "
$ws.Range("E11").Value = "  +3.87%  "
$ws.Range("D12").Value = "0.0000341"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "42.07"
$ws.Range("E13").Value = "  +8.54%  "
$ws.Range("D14").Value = "10.45"
$ws.Range("E14").Value = "  +13.40%  "
$ws.Range("D15").Value = "4.448.76"
$ws.Range("E15").Value = "  +10.10%  "
$ws.Range("D16").Value = "15.82"
$ws.Range("E16").Value = "  +26.13%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.847.41"
$ws.Range("E17").Value = "  +9.83%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.138"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "20.06"
$ws.Range("E19").Value = "  +7.07%  "
$ws.Range("E20").Value = "  +8.18%  "
$ws.Range("D21").Value = "66.680.15"
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").Value = "415.11"
$ws.Range("E22").Value = "  +5.68%  "
$ws.Range("D23").Value = "15.17"
$ws.Range("E23").Value = "  +8.98%  "
$ws.Range("D24").Value = "84.95"
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("D25").Value = "3.10"
$ws.Range("E25").Value = "  +8.68%  "
$ws.Range("E26").Value = "  +13.43%  "
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +14.77%  "
$ws.Range("D28").Value = "3.27"
$ws.Range("E28").Value = "  +9.27%  "
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "9.26"
$ws.Range("E30").Value = "  +35.24%  "
$ws.Range("D31").Value = "719.68"
$ws.Range("E31").Value = "  +7.79%  "
$ws.Range("D32").Value = "13.88"
$ws.Range("E32").Value = "  +15.47%  "
$ws.Range("E33").Value = "  +14.47%  "
$ws.Range("D34").Value = "2.80"
$ws.Range("E34").Value = "  +6.82%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.81"
$ws.Range("E36").Value = "  +44.13%  "
$ws.Range("D37").Value = "39.20"
$ws.Range("E37").Value = "  +6.63%  "
$ws.Range("D38").Value = "0.153"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "55.74"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").Value = "0.0₃0750"
$ws.Range("E40").Value = "  +17.85%  "
$ws.Range("E41").Value = "  +6.23%  "
$ws.Range("D42").Value = "2.92"
$ws.Range("E42").Value = "  +7.10%  "
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  +6.88%  "
$ws.Range("D45").Value = "0.137"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("E46").Value = "  +10.25%  "
$ws.Range("D47").Value = "0.322"
$ws.Range("E47").Value = "  +16.20%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.86"
$ws.Range("E48").Value = "  +5.76%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "142.60"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "2.62"
$ws.Range("E50").Value = "  +4.57%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "2.05"
$ws.Range("E51").Value = "  +5.61%  "
